$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")
$ws.Activate()

# Insert a new (blank) column before column N (14) - this is what Excel
# does when a user right-clicks a column header and chooses "Insert".
# Everything from N onward (N..P) shifts right by one (to O..Q).
$ws.Columns("N").Insert()

# New column inherits the width of the column to its left (M)
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Move the selection to T7 as recorded after the edit
$ws.Range("T7").Select()
